$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.016154050827026
$ws.Range("B1").Value = 2.483772993087769
$ws.Range("C1").Value = 2.584987878799438
$ws.Range("D1").Value = 3.293428182601929
$ws.Range("E1").Value = 1.080365061759949
